$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns.Item(3).Insert()

Write-Output $ws.UsedRange.Rows.Count
Write-Output $ws.UsedRange.Columns.Count
